# "Generate Report for Handback"
# The localization handback finished: the status moves from "Ready for
# handoff" to "Handed back: in sync with en-US", the per-language sheets
# get their "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns populated (with a hyperlink on the
# target-file cell), and a few report columns are widened so the longer
# values are readable.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
$mdFileName = "484f6901-4581-47a0-b030-22d7ae1b9968.md"
$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3e479601223c5b8a61140f38f253c5b91a12da47/e2e/484f6901-4581-47a0-b030-22d7ae1b9968.md"

# ColumnWidth (character units) round-trips through Excel's pixel snapping,
# which always adds ~5/6 of a character to whatever is assigned. Back that
# offset out so the stored width lands as close as possible to the target.
function Set-ColumnWidthPrecise($col, $targetWidth) {
    $col.ColumnWidth = $targetWidth - (5 / 6)
}

# ---------------------------------------------------------------------
# Overview sheet: widen the two per-language status columns (E, F)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
Set-ColumnWidthPrecise $wsOverview.Columns.Item(5) 29.9777047293527
Set-ColumnWidthPrecise $wsOverview.Columns.Item(6) 29.9777047293527

# ---------------------------------------------------------------------
# zh-cn sheet: handback completed 2016-08-31 13:14:39
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $statusText
$wsZhCn.Range("J2").Value = "484f6901-4581-47a0-b030-22d7ae1b9968.78ea64c13f558f3426926f30bb3dfe80916f2890.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-31 13:14:39"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $mdUrl, "", "", $mdFileName)
Set-ColumnWidthPrecise $wsZhCn.Columns.Item(3) 29.9777047293527
Set-ColumnWidthPrecise $wsZhCn.Columns.Item(9) 40
Set-ColumnWidthPrecise $wsZhCn.Columns.Item(10) 40

# ---------------------------------------------------------------------
# de-de sheet: handback completed 2016-08-31 13:14:47
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $statusText
$wsDeDe.Range("J2").Value = "484f6901-4581-47a0-b030-22d7ae1b9968.78ea64c13f558f3426926f30bb3dfe80916f2890.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-31 13:14:47"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $mdUrl, "", "", $mdFileName)
Set-ColumnWidthPrecise $wsDeDe.Columns.Item(3) 29.9777047293527
Set-ColumnWidthPrecise $wsDeDe.Columns.Item(9) 40
Set-ColumnWidthPrecise $wsDeDe.Columns.Item(10) 40
